$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest cryptos snapshot

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "27.390.18"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -1.09%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.710.51"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -1.36%  "

$ws.Range("E4").Value = "  +0.03%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "224.39"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -1.37%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.5341"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -2.08%  "

$ws.Range("E7").Value = "  +0.13%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.2681"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -2.23%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.06618"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -1.41%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "21.00"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -3.87%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07614"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -2.05%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "4.558"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -2.95%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.703.07"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -2.04%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "1.945.77"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.42%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.5781"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -3.37%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.0₅8188"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -2.78%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "67.87"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -1.81%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "27.357.47"
$cell.Style = "Normal"

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "217.29"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -4.07%  "

$ws.Range("E20").Value = "  +0.12%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "4.677"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -3.13%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "10.47"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -3.93%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "5.976"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -3.94%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "1.005"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.07%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "142.35"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -3.92%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "1.741"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +1.03%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.1215"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -2.87%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "7.272"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -2.68%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "16.29"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -4.73%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.05405"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -5.28%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.292"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -1.56%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.503"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -5.39%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "3.432"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -2.51%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.646"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -2.43%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "2.879"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +0.85%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.9519"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -2.41%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.416"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.84%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.5872"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -1.93%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.01637"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -1.93%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "5.859"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -1.26%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.047.13"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -0.18%  "

$ws.Range("E42").Value = "  +0.11%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.8421"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -0.97%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "101.13"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.51%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.852.86"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -1.35%  "

$ws.Range("E46").Value = "  +1.83%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "58.06"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -2.75%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.4516"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +1.99%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.007"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.87%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "8.080"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -2.50%  "

$ws.Range("E51").Value = "  -1.85%  "

